# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and refresh the
# "Latest Handoff Datetime" timestamps on each localization sheet, then
# widen the Status columns to fit the longer text.

$wb = $excel.ActiveWorkbook

# New, wider "Status" column width (character units) that reproduces the
# target stored column width of the new, longer "Ready for handoff" text.
$statusColWidth = 16.33

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-13 08:45:50"
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-13 08:45:43"
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-13 08:45:50"
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
